# RC-DE.schema.docx auto-generation update.
#
# The "DistributionElement" table is trimmed back down to just its header
# row, and the whole "sender" / "recipient" sub-sections (Heading1
# paragraph + following 2-row reference table) are removed entirely.

$d = $word.ActiveDocument

# --- 1. Drop the "recipient" table (3rd table) and the "sender" table
#        (2nd table) in full. Go from the end backwards so earlier table
#        indices stay valid while we work. ---
$d.Tables.Item(3).Delete()
$d.Tables.Item(2).Delete()

# --- 2. Trim the first ("DistributionElement") table down to just its
#        header row by repeatedly deleting the 2nd row until only the
#        header remains. ---
$t1 = $d.Tables.Item(1)
while ($t1.Rows.Count -gt 1) {
    $t1.Rows.Item(2).Delete()
}

# --- 3. Remove the now-orphaned "sender" and "recipient" Heading1
#        paragraphs that used to introduce the deleted tables.
#        ($d.Content.Paragraphs is used instead of $d.Paragraphs so the
#        collection is re-walked fresh after the structural edits above.)
#        Matched on style + "outside any table" as well as exact text so
#        this can't accidentally catch the like-named table-cell values
#        (e.g. the "sender"/"recipient" field-name cells still present in
#        the trimmed first table's body rows, if any were ever to remain).
for ($i = $d.Content.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Content.Paragraphs.Item($i)
    $txt = $p.Range.Text.TrimEnd([char]13, [char]7)
    $isHeading1 = ($p.Style.NameLocal -eq "Heading 1")
    $inTable = $p.Range.Information(12)
    if (($txt -eq "sender" -or $txt -eq "recipient") -and $isHeading1 -and -not $inTable) {
        $p.Range.Delete()
    }
}
